$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "CollectionAgency"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "bnb"

# --- Update header/data cells on the CollectionAgency sheet ---
# (write order chosen so new shared strings land in the same append order
#  as the target workbook)
$ws1.Range("B2").Value = "TC_ProcessSheet_CollectionAgency"
$ws1.Range("C2").Value = "Role5"
$ws1.Range("C1").Value = "Role Name"

# --- Cosmetic view-state changes ---
# Widen column C on CollectionAgency sheet
$ws1.Columns.Item(3).ColumnWidth = 14.3

# Update the selection on the "bnb" sheet first (without leaving it active)
$ws2.Range("B20").Select()

# Re-activate CollectionAgency and update its selection / scroll position
$ws1.Activate()
$ws1.Range("A9").Select()
